$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New entry: Corona Lüftung Aufgaben, 26.10.2020, 11:45 - 17:45
$ws.Range("A6").Value = 44130
$ws.Range("B6").Value = 0.48958333333333331
$ws.Range("C6").Value = 0.73958333333333337
$ws.Range("D6").Formula = "=HOUR(C6)+MINUTE(C6)/60-HOUR(B6)-MINUTE(B6)/60+D5"
$ws.Range("E6").Value = "Corona Lüftung Aufgaben"

# Match the number formats used by the row above (reuse existing styles)
$ws.Range("A5").Copy()
$ws.Range("A6").PasteSpecial(-4122)

$ws.Range("B5:C5").Copy()
$ws.Range("B6:C6").PasteSpecial(-4122)

$ws.Range("D5").Copy()
$ws.Range("D6").PasteSpecial(-4122)

$ws.Range("E11").Select()
